# Added test cases for Recurring billing. Maintenance of other test cases.
#
# The "RMA Details Maintenance Grid" sheet keeps the most recent RMA test
# case identifiers in columns E (RMA number), F (RMA line number) and J
# (Salesforce record id) for rows 2-4. A new RMA test case (O9J6) replaces
# the previous one (MO6T) that was there before.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Row 2 - line 1
$ws.Range("E2").Value = "RMA-O9J6-001"
$ws.Range("F2").Value = "RMA-O9J6-1-1"
$ws.Range("J2").Value = "a7s5f000000xKBvAAM"

# Row 3 - line 2
$ws.Range("E3").Value = "RMA-O9J6-002"
$ws.Range("F3").Value = "RMA-O9J6-1-2"
$ws.Range("J3").Value = "a7s5f000000xKBwAAM"

# Row 4 - line 3
$ws.Range("E4").Value = "RMA-O9J6-003"
$ws.Range("F4").Value = "RMA-O9J6-1-3"
$ws.Range("J4").Value = "a7s5f000000xKBxAAM"
